$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Build the "template" sheet by copying the existing booking-form sheet
#    ("Sheet2") to a new tab placed before "Sheet1". This new tab becomes
#    "Sheet5" - it starts as an exact duplicate (values + styles + widths)
#    of Sheet2's A1:H2 booking form.
# ---------------------------------------------------------------------------
$src = $wb.Worksheets.Item("Sheet2")
$beforeTarget = $wb.Worksheets.Item("Sheet1")
$src.Copy($beforeTarget)
$sheet5 = $wb.Worksheets.Item("Sheet2 (2)")
$sheet5.Name = "Sheet5"

# Update the check-in / check-out dates on the template to the new values
# used by the booking-with-occupancy sheets (5-Aug-2025 / 5-Sep-2025).
$sheet5.Range("G2").Value = 45874
$sheet5.Range("H2").Value = 45905

# Add the new "adults / children per room" columns. The write order below
# matters: it fixes the order new shared strings are interned in, matching
# the order the workbook's sharedStrings table was built in (header cells
# across I/J, then the "expectedFinalPrice" header that was later moved to
# column L, then the row-2 occupancy values, then the "expectedTotalPrice"
# header and its numeric values in K/L).
$sheet5.Range("I1").Value = "adultsPerRoom"
$sheet5.Range("J1").Value = "childrenPerRoom"
$sheet5.Range("L1").Value = "expectedFinalPrice"
$sheet5.Range("I2").Value = "1 - One"
$sheet5.Range("J2").Value = "0 - None"
$sheet5.Range("K1").Value = "expectedTotalPrice"
$sheet5.Range("K2").Value = 125
$sheet5.Range("L2").Value = 135

# ---------------------------------------------------------------------------
# 2. Duplicate the fully built template (12 columns, A:L) twice more:
#      - "Sheet3" sits right after "Sheet5", ahead of the original sheets.
#      - "Sheet4" is appended after the original "Sheet2" tab and is the
#        only one of the three that keeps the expectedTotalPrice /
#        expectedFinalPrice columns (K:L).
# ---------------------------------------------------------------------------
$beforeTarget2 = $wb.Worksheets.Item("Sheet1")
$sheet5.Copy($beforeTarget2)
$sheet3 = $wb.Worksheets.Item("Sheet5 (2)")
$sheet3.Name = "Sheet3"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet5b = $wb.Worksheets.Item("Sheet5")
$sheet5b.Copy($null, $lastSheet)
$sheet4 = $wb.Worksheets.Item("Sheet5 (2)")
$sheet4.Name = "Sheet4"

# ---------------------------------------------------------------------------
# 3. "Sheet5" and "Sheet3" only go up to column J (noOfRooms .. childrenPerRoom)
#    - trim the expectedTotalPrice / expectedFinalPrice columns back off of
#    them; "Sheet4" is left with all 12 columns.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Sheet5").Range("K1:L2").Clear()
$wb.Worksheets.Item("Sheet3").Range("K1:L2").Clear()

# ---------------------------------------------------------------------------
# 4. View/selection bookkeeping so the final file matches the saved state:
#    Sheet5 is the active tab; each new sheet keeps the selection where the
#    user last left it.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Sheet3").Range("J2").Select()
$wb.Worksheets.Item("Sheet4").Range("K7").Select()

$sheet5Final = $wb.Worksheets.Item("Sheet5")
$sheet5Final.Activate()
$sheet5Final.Range("K1:L1048576").Select()

Write-Host "Final sheet order:"
foreach ($s in $wb.Worksheets) { Write-Host $s.Name }
